$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh of the scraped cryptocurrency "Price" (D) and "Volume(1h)" (E)
# columns to the latest snapshot pulled by the scheduled scraper run.
# Prices / percentages are plain text cells (t="inlineStr"): a leading
# apostrophe is used only where the new text would otherwise parse as a
# plain number (e.g. '1.00', '0.0000226'), exactly like typing it into Excel,
# so the cell keeps its original text type instead of becoming numeric.

$ws.Range("D2").Value = "60.949.57"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "2.919.93"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'593.74"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").Value = "'145.83"
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.506"
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("D9").Value = "'6.84"
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("E11").Value = "  -1.95%  "
$ws.Range("D12").Value = "'0.0000226"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "'33.63"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").Value = "3.403.47"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").Value = "60.930.48"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "'6.70"
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").Value = "2.918.40"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").Value = "'430.85"
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("D20").Value = "'13.39"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").Value = "'10.93"
$ws.Range("E24").Value = "  -1.51%  "
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("D26").Value = "'11.96"
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  +5.06%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "'2.62"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").Value = "'7.06"
$ws.Range("E31").Value = "  -2.19%  "
$ws.Range("D32").Value = "'26.42"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").Value = "0.0₃0852"
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").Value = "'5.63"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("E37").Value = "  +1.82%  "
$ws.Range("D38").Value = "'1.99"
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("D39").Value = "'0.122"
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("E40").Value = "  -1.65%  "
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("D42").Value = "'39.97"
$ws.Range("E42").Value = "  -3.48%  "
$ws.Range("D43").Value = "'377.11"
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("D45").Value = "2.700.43"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("D46").Value = "'131.23"
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("D48").Value = "'23.99"
$ws.Range("E48").Value = "  -5.61%  "
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("E50").Value = "  -3.75%  "
$ws.Range("E51").Value = "  +1.61%  "
